$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark that used to sit right after the
#    "... (para tipos java)." bullet.
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------------
# 2. Append two new bulleted list items right after the "... minutos y
#    segundos." paragraph (and before the "Datatypes:" paragraph), and move
#    the "_GoBack" bookmark to the end of the new text, matching what Word
#    does when the cursor is left there after the last edit.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$targetIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -match "minutos y segundos\.") {
        $targetIndex = $i
    }
}

$tr = $d.Paragraphs.Item($targetIndex).Range
$tr.InsertParagraphAfter()

$newRange = $d.Paragraphs.Item($targetIndex + 1).Range
$newRange.Collapse(1)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p>' +
    '<w:pPr>' +
      '<w:pStyle w:val="Prrafodelista"/>' +
      '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' +
      '<w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>' +
    '</w:pPr>' +
    '<w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Cambiar array por </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Collection</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">en </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>collecciones</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> Java.</w:t></w:r>' +
  '</w:p>' +
  '<w:p>' +
    '<w:pPr>' +
      '<w:pStyle w:val="Prrafodelista"/>' +
      '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' +
      '<w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>' +
    '</w:pPr>' +
    '<w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>TODAS COLLECCIONES NOTNULL.</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>' +
  '</w:p>' +
  '</w:body>' +
  '</w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

[void]$newRange.InsertXML($xml)
